# Insert a new weekly price record at row 155 for
# "Terminal La Palmera de La Serena - Zanahoria", pushing the existing
# rows 155:224 down to 156:225 (the sheet's dimension grows from
# A1:R224 to A1:R225).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 155; this shifts rows
# 155-224 down to 156-225 and stretches the used range accordingly.
$ws.Rows.Item(155).Insert()

# Populate the newly inserted row 155 with the new record.
$ws.Range("A155").Value = 8
$ws.Range("B155").Value = "Terminal La Palmera de La Serena"
$ws.Range("C155").Value = "Coquimbo"
$ws.Range("D155").Value = 44510
$ws.Range("E155").Value = 4
$ws.Range("F155").Value = 100114013
$ws.Range("G155").Value = "Zanahoria"
$ws.Range("H155").Value = "Sin especificar"
$ws.Range("I155").Value = "Primera"
$ws.Range("J155").Value = 800
$ws.Range("K155").Value = 6500
$ws.Range("L155").Value = 7000
$ws.Range("M155").Value = 6750
$ws.Range("N155").Value = "`$/saco 20 kilos"
$ws.Range("O155").Value = "Provincia del Elquí"
$ws.Range("P155").Value = 338
$ws.Range("Q155").Value = 20
$ws.Range("R155").Value = "Hortaliza"
